$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.0995906397018391"
$ws.Range("C2").Value = [double]"0.968106555874626"
$ws.Range("D2").Value = [double]"0.000488788415714548"
$ws.Range("E2").Value = [double]"0.993890144803568"
$ws.Range("F2").Value = [double]"0.00837050161911163"
$ws.Range("G2").Value = [double]"0.955398057066048"
$ws.Range("H2").Value = [double]"0.0231563511944767"
$ws.Range("I2").Value = [double]"0.276348750534612"
$ws.Range("J2").Value = [double]"0.000244394207857274"
$ws.Range("K2").Value = [double]"0.95246532657176"
$ws.Range("L2").Value = [double]"0.991996089692674"
$ws.Range("M2").Value = [double]"0.874198081505468"
$ws.Range("N2").Value = [double]"0.943606036536934"
$ws.Range("O2").Value = [double]"0.000855379727500458"
$ws.Range("P2").Value = [double]"0.999205718824464"
$ws.Range("Q2").Value = [double]"6.10985519643184e-05"
$ws.Range("R2").Value = [double]"0.999022423168571"
$ws.Range("S2").Value = [double]"6.10985519643184e-05"
$ws.Range("T2").Value = [double]"0.00403250442964502"
$ws.Range("U2").Value = [double]"0.995295411498747"
$ws.Range("V2").Value = [double]"0.988024683814994"
$ws.Range("W2").Value = [double]"0.0887761960041547"
$ws.Range("X2").Value = [double]"0.000733182623571821"

$ws.Range("B3").Value = [double]"0.0530946416569927"
$ws.Range("C3").Value = [double]"0.000305492759821592"
$ws.Range("D3").Value = [double]"0.00122197103928637"
$ws.Range("E3").Value = [double]"0.00128306959125069"
$ws.Range("F3").Value = [double]"0.00775951609946844"
$ws.Range("G3").Value = [double]"0.0152135394391153"
$ws.Range("H3").Value = [double]"0.0389197776012708"
$ws.Range("I3").Value = [double]"6.10985519643184e-05"
$ws.Range("J3").Value = [double]"0.000549886967678866"
$ws.Range("K3").Value = [double]"0.000916478279464777"
$ws.Range("L3").Value = [double]"0.00378811022178774"
$ws.Range("M3").Value = [double]"0.00409360298160934"
$ws.Range("N3").Value = [double]"0.000610985519643184"
$ws.Range("P3").Value = [double]"6.10985519643184e-05"
$ws.Range("Q3").Value = [double]"0.917883546159956"
$ws.Range("R3").Value = [double]"6.10985519643184e-05"
$ws.Range("S3").Value = [double]"0.999450113032321"
$ws.Range("T3").Value = [double]"0.985702938840349"
$ws.Range("U3").Value = [double]"0.00342151891000183"
$ws.Range("V3").Value = [double]"0.000366591311785911"
$ws.Range("W3").Value = [double]"0.000183295655892955"
$ws.Range("X3").Value = [double]"0.000122197103928637"

$ws.Range("B4").Value = [double]"0.774668540355594"
$ws.Range("C4").Value = [double]"0.027433249831979"
$ws.Range("D4").Value = [double]"0.000549886967678866"
$ws.Range("E4").Value = [double]"0.00354371601393047"
$ws.Range("F4").Value = [double]"0.982709109794098"
$ws.Range("G4").Value = [double]"0.00879819148286186"
$ws.Range("H4").Value = [double]"0.547381927048329"
$ws.Range("I4").Value = [double]"0.717358098613063"
$ws.Range("J4").Value = [double]"0.000244394207857274"
$ws.Range("K4").Value = [double]"0.0430133805828802"
$ws.Range("L4").Value = [double]"0.00348261746196615"
$ws.Range("M4").Value = [double]"0.117309219771491"
$ws.Range("N4").Value = [double]"0.0543777112482434"
$ws.Range("O4").Value = [double]"0.998961324616607"
$ws.Range("P4").Value = [double]"0.000733182623571821"
$ws.Range("R4").Value = [double]"0.000672084071607503"
$ws.Range("S4").Value = [double]"6.10985519643184e-05"
$ws.Range("T4").Value = [double]"6.10985519643184e-05"
$ws.Range("U4").Value = [double]"0.000733182623571821"
$ws.Range("V4").Value = [double]"0.0106922465937557"
$ws.Range("W4").Value = [double]"0.905724934319057"
$ws.Range("X4").Value = [double]"0.9991446202725"

$ws.Range("B5").Value = [double]"0.068308181096108"
$ws.Range("C5").Value = [double]"0.00311602615018024"
$ws.Range("D5").Value = [double]"0.997678255025356"
$ws.Range("E5").Value = [double]"0.000366591311785911"
$ws.Range("F5").Value = [double]"0.000183295655892955"
$ws.Range("G5").Value = [double]"0.0197348322844749"
$ws.Range("H5").Value = [double]"0.384065497647706"
$ws.Range("I5").Value = [double]"0.000427689863750229"
$ws.Range("J5").Value = [double]"0.998900226064642"
$ws.Range("K5").Value = [double]"0.00183295655892955"
$ws.Range("L5").Value = [double]"0.000183295655892955"
$ws.Range("M5").Value = [double]"0.000122197103928637"
$ws.Range("N5").Value = [double]"6.10985519643184e-05"
$ws.Range("O5").Value = [double]"6.10985519643184e-05"
$ws.Range("Q5").Value = [double]"0.0799169059693285"
$ws.Range("R5").Value = [double]"6.10985519643184e-05"
$ws.Range("S5").Value = [double]"0.000366591311785911"
$ws.Range("T5").Value = [double]"0.00947027555446936"
$ws.Range("U5").Value = [double]"0.000122197103928637"
$ws.Range("V5").Value = [double]"0.000427689863750229"
$ws.Range("W5").Value = [double]"0.00201625221482251"
